$d = $word.ActiveDocument
$cursor = 0
$pairs = @(
    @("76-65=11", "19-16=3"),
    @("57-3=54", "68-58=10"),
    @("66-57=9", "41+33=74"),
    @("20+25=45", "86-85=1"),
    @("36+9=45", "15+15=30"),
    @("21-2=19", "23+73=96"),
    @("43+7=50", "71-57=14"),
    @("25+5=30", "93-75=18"),
    @("76-60=16", "38-17=21"),
    @("62+15=77", "1+35=36"),
    @("21+55=76", "75-62=13"),
    @("89-60=29", "43+28=71"),
    @("27+45=72", "22+14=36"),
    @("7+5=12", "14+84=98"),
    @("31+21=52", "31+48=79"),
    @("87+7=94", "9+20=29"),
    @("91-56=35", "73-46=27"),
    @("85-56=29", "86+4=90"),
    @("77+2=79", "69-44=25"),
    @("42-2=40", "56-40=16"),
    @("89-38=51", "61+0=61"),
    @("99-76=23", "43-22=21"),
    @("27-16=11", "28+15=43"),
    @("78-29=49", "93-66=27"),
    @("75-64=11", "7+59=66"),
    @("76-62=14", "15+24=39"),
    @("94-83=11", "72-28=44"),
    @("71-54=17", "56-23=33"),
    @("98-64=34", "16+66=82"),
    @("40+40=80", "49+12=61"),
    @("95-25=70", "15-6=9"),
    @("33+18=51", "21+48=69"),
    @("57+41=98", "87-67=20"),
    @("82-70=12", "74-22=52"),
    @("40+6=46", "18+62=80"),
    @("34-25=9", "59+4=63"),
    @("57-4=53", "22+62=84"),
    @("85-51=34", "46+0=46"),
    @("57+5=62", "32+57=89"),
    @("22-7=15", "61+19=80"),
    @("33+18=51", "63-34=29"),
    @("47-26=21", "89-12=77"),
    @("92-25=67", "99-65=34"),
    @("90-69=21", "69+12=81"),
    @("52+19=71", "8+43=51"),
    @("0+84=84", "45+22=67"),
    @("48+8=56", "67-44=23"),
    @("60+23=83", "23-6=17"),
    @("43-23=20", "44+9=53"),
    @("57-31=26", "6+56=62"),
    @("6+63=69", "53+15=68"),
    @("8+37=45", "20+11=31"),
    @("96-23=73", "68-50=18"),
    @("30+41=71", "85+2=87"),
    @("77-49=28", "45+17=62"),
    @("89-78=11", "79+16=95"),
    @("12+1=13", "56-26=30"),
    @("77-60=17", "82-14=68"),
    @("49+1=50", "36+18=54"),
    @("41+21=62", "69-34=35"),
    @("58-3=55", "9+24=33"),
    @("86-50=36", "99-31=68"),
    @("30+25=55", "76-27=49"),
    @("8+59=67", "83+7=90"),
    @("65-14=51", "67-51=16"),
    @("60+22=82", "67+12=79"),
    @("24+41=65", "55+18=73"),
    @("95-27=68", "56-56=0"),
    @("8+22=30", "19-4=15"),
    @("17+61=78", "52+43=95"),
    @("83-47=36", "57+29=86"),
    @("94-1=93", "19-9=10"),
    @("89-23=66", "71-25=46"),
    @("58-7=51", "78-36=42"),
    @("80-46=34", "24+72=96"),
    @("81-57=24", "28+29=57"),
    @("4+94=98", "98-37=61"),
    @("64+28=92", "90+2=92"),
    @("38+46=84", "8+24=32"),
    @("55-2=53", "54-37=17"),
    @("18+31=49", "49+24=73"),
    @("68-26=42", "21+19=40"),
    @("38-36=2", "22-2=20"),
    @("2+44=46", "24+8=32"),
    @("8+6=14", "65+24=89"),
    @("43-39=4", "13+63=76"),
    @("10+41=51", "0+36=36"),
    @("95-63=32", "22+50=72"),
    @("45+18=63", "44-5=39"),
    @("91-57=34", "96-87=9"),
    @("13+66=79", "5+36=41"),
    @("82-45=37", "71-26=45"),
    @("93-12=81", "33+45=78"),
    @("75-31=44", "25+68=93"),
    @("67+19=86", "48+14=62"),
    @("81-21=60", "92-72=20"),
    @("58+27=85", "78-65=13"),
    @("85+0=85", "8+10=18"),
    @("87-85=2", "6-4=2"),
    @("89-61=28", "92-70=22")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Range($cursor, $d.Content.End)
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $new
        $cursor = $rng.End
    } else {
        Write-Host "NOT FOUND: $old"
    }
}
